# Apply targeted cell updates to rows 26-35 per the source diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 26
$ws.Range("A26").Value = 112076815
$ws.Range("B26").Value = 96735
$ws.Range("D26").Value = "VU"
$ws.Range("E26").Value = 220787
$ws.Range("F26").Value = "Knärot"
$ws.Range("G26").Value = "Goodyera repens"
$ws.Range("H26").Value = "(L.) R. Br."
$ws.Range("Q26").Value = 485636
$ws.Range("R26").Value = 7005629

# Row 27
$ws.Range("B27").Value = 86371

# Row 28
$ws.Range("A28").Value = 112076818
$ws.Range("B28").Value = 89047
$ws.Range("D28").Value = "NT"
$ws.Range("E28").Value = 3286
$ws.Range("F28").Value = "Flattoppad klubbsvamp"
$ws.Range("G28").Value = "Clavariadelphus truncatus"
$ws.Range("H28").Value = "(Quél.) Donk"
$ws.Range("Q28").Value = 485662
$ws.Range("R28").Value = 7005637

# Row 29
$ws.Range("A29").Value = 112076816
$ws.Range("B29").Value = 96735
$ws.Range("D29").Value = "VU"
$ws.Range("E29").Value = 220787
$ws.Range("F29").Value = "Knärot"
$ws.Range("G29").Value = "Goodyera repens"
$ws.Range("H29").Value = "(L.) R. Br."
$ws.Range("Q29").Value = 485618
$ws.Range("R29").Value = 7005614

# Row 30
$ws.Range("A30").Value = 112076820
$ws.Range("B30").Value = 98980
$ws.Range("D30").Value = "LC"
$ws.Range("E30").Value = 222498
$ws.Range("F30").Value = "Blåsippa"
$ws.Range("G30").Value = "Hepatica nobilis"
$ws.Range("H30").Value = "Schreb."
$ws.Range("Q30").Value = 485536
$ws.Range("R30").Value = 7005851

# Row 31
$ws.Range("A31").Value = 112076814
$ws.Range("B31").Value = 90235
$ws.Range("D31").Value = "LC"
$ws.Range("E31").Value = 3298
$ws.Range("F31").Value = "Trådticka"
$ws.Range("G31").Value = "Climacocystis borealis"
$ws.Range("H31").Value = "(Fr.) Kotl. & Pouzar"
$ws.Range("Q31").Value = 485714
$ws.Range("R31").Value = 7005798

# Row 32
$ws.Range("A32").Value = 112076817
$ws.Range("B32").Value = 96735
$ws.Range("Q32").Value = 485596
$ws.Range("R32").Value = 7005613

# Row 33
$ws.Range("A33").Value = 112076813
$ws.Range("B33").Value = 89553
$ws.Range("E33").Value = 1202
$ws.Range("F33").Value = "Ullticka"
$ws.Range("G33").Value = "Phellinidium ferrugineofuscum"
$ws.Range("H33").Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Range("Q33").Value = 485752
$ws.Range("R33").Value = 7005707

# Row 34
$ws.Range("B34").Value = 86371

# Row 35
$ws.Range("B35").Value = 98980
